$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 177.76923
$ws.Range("I33").Value = 95.3
$ws.Range("K33").Value = 95.3
$ws.Range("M33").Value = 133.7
$ws.Range("H137").Value = 1460.0834
$ws.Range("I137").Value = 961.8
$ws.Range("J137").Value = 1816
$ws.Range("K137").Value = 2885.4
$ws.Range("L137").Value = 5448
$ws.Range("M137").Value = -335.3999999999996
$ws.Range("N137").Value = -10548
$ws.Range("H138").Value = 1481.95
$ws.Range("I138").Value = 668.5714
$ws.Range("J138").Value = 1919.9231
$ws.Range("K138").Value = 2005.7142
$ws.Range("L138").Value = 5759.7693
$ws.Range("M138").Value = 3134.2858
$ws.Range("N138").Value = -16039.7693

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4045.8667
$ws.Range("I32").Value = 3677.673
$ws.Range("J32").Value = 6439.125
$ws.Range("K32").Value = 3677.673
$ws.Range("L32").Value = 6439.125
$ws.Range("M32").Value = -3390.673
$ws.Range("N32").Value = -7013.125
$ws.Range("H61").Value = 1276.5555
$ws.Range("I61").Value = 1122.0667
$ws.Range("J61").Value = 2049
$ws.Range("K61").Value = 1122.0667
$ws.Range("L61").Value = 2049
$ws.Range("M61").Value = -910.0667000000001
$ws.Range("N61").Value = -2473
$ws.Range("H74").Value = 1964.8
$ws.Range("I74").Value = 1012
$ws.Range("J74").Value = 2600
$ws.Range("K74").Value = 1012
$ws.Range("L74").Value = 2600
$ws.Range("M74").Value = -138
$ws.Range("N74").Value = -4348
$ws.Range("H77").Value = 1964.8
$ws.Range("I77").Value = 1012
$ws.Range("J77").Value = 2600
$ws.Range("K77").Value = 5060
$ws.Range("L77").Value = 13000
$ws.Range("M77").Value = -692
$ws.Range("N77").Value = -21736
$ws.Range("H132").Value = 2239.7778
$ws.Range("I132").Value = 1879.8572
$ws.Range("K132").Value = 5639.571599999999
$ws.Range("M132").Value = -3109.571599999999
$ws.Range("H136").Value = 1276.5555
$ws.Range("I136").Value = 1122.0667
$ws.Range("J136").Value = 2049
$ws.Range("K136").Value = 3366.2001
$ws.Range("L136").Value = 6147
$ws.Range("M136").Value = -816.2001
$ws.Range("N136").Value = -11247

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1154.1852
$ws.Range("I31").Value = 1153.3208
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 1153.3208
$ws.Range("L31").Value = 1200
$ws.Range("M31").Value = -858.3208
$ws.Range("N31").Value = -1790
$ws.Range("H34").Value = 1154.1852
$ws.Range("I34").Value = 1153.3208
$ws.Range("J34").Value = 1200
$ws.Range("K34").Value = 1153.3208
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = -951.3208
$ws.Range("N34").Value = -1604
$ws.Range("H94").Value = 1000.8182
$ws.Range("I94").Value = 787.4
$ws.Range("J94").Value = 1178.6666
$ws.Range("K94").Value = 787.4
$ws.Range("L94").Value = 1178.6666
$ws.Range("M94").Value = -336.4
$ws.Range("N94").Value = -2080.6666
$ws.Range("H132").Value = 1750.25
$ws.Range("I132").Value = 1296.6666
$ws.Range("K132").Value = 3889.9998
$ws.Range("M132").Value = -1359.9998
$ws.Range("H134").Value = 1025.9131
$ws.Range("I134").Value = 864.2857
$ws.Range("J134").Value = 1277.3334
$ws.Range("K134").Value = 2592.8571
$ws.Range("L134").Value = 3832.0002
$ws.Range("M134").Value = -57.85710000000017
$ws.Range("N134").Value = -8902.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 975612.25
$ws.Range("I4").Value = 224804.25
$ws.Range("J4").Value = 1351016.2
$ws.Range("K4").Value = 674412.75
$ws.Range("L4").Value = 4053048.6
$ws.Range("M4").Value = -674300.75
$ws.Range("N4").Value = -4053272.6
$ws.Range("H5").Value = 1191.258
$ws.Range("I5").Value = 1289.7693
$ws.Range("J5").Value = 679
$ws.Range("K5").Value = 3869.3079
$ws.Range("L5").Value = 2037
$ws.Range("M5").Value = -3757.3079
$ws.Range("N5").Value = -2261
$ws.Range("H34").Value = 2232
$ws.Range("I34").Value = 1647
$ws.Range("J34").Value = 2700
$ws.Range("K34").Value = 4941
$ws.Range("L34").Value = 8100
$ws.Range("M34").Value = -4857
$ws.Range("N34").Value = -8268
$ws.Range("H125").Value = 998
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 18184666
$ws.Range("I131").Value = 200000370
$ws.Range("J131").Value = 3095.54
$ws.Range("K131").Value = 600001110
$ws.Range("L131").Value = 9286.619999999999
$ws.Range("M131").Value = -599996070
$ws.Range("N131").Value = -19366.62
$ws.Range("H135").Value = 1191.258
$ws.Range("I135").Value = 1289.7693
$ws.Range("J135").Value = 679
$ws.Range("K135").Value = 11607.9237
$ws.Range("L135").Value = 6111
$ws.Range("M135").Value = -9072.923699999999
$ws.Range("N135").Value = -11181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 25489.5
$ws.Range("J52").Value = 25980
$ws.Range("L52").Value = 25980
$ws.Range("N52").Value = -26498
$ws.Range("H80").Value = 2739
$ws.Range("I80").Value = 2368.889
$ws.Range("J80").Value = 3041.818
$ws.Range("K80").Value = 2368.889
$ws.Range("L80").Value = 3041.818
$ws.Range("M80").Value = -1370.889
$ws.Range("N80").Value = -5037.818
$ws.Range("H83").Value = 2739
$ws.Range("I83").Value = 2368.889
$ws.Range("J83").Value = 3041.818
$ws.Range("K83").Value = 11844.445
$ws.Range("L83").Value = 15209.09
$ws.Range("M83").Value = -6852.445
$ws.Range("N83").Value = -25193.09
$ws.Range("H107").Value = 740312.5
$ws.Range("I107").Value = 1069057
$ws.Range("K107").Value = 1069057
$ws.Range("M107").Value = -1067137
$ws.Range("H135").Value = 29874.625
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 29874.625
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 29874.625
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -40014.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 227322.36
$ws.Range("J2").Value = 97475.86
$ws.Range("L2").Value = 97475.86
$ws.Range("N2").Value = -97699.86
$ws.Range("H55").Value = 228.65517
$ws.Range("I55").Value = 189.4375
$ws.Range("K55").Value = 189.4375
$ws.Range("M55").Value = -16.4375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H92").Value = 14337.25
$ws.Range("J92").Value = 14337.25
$ws.Range("L92").Value = 14337.25
$ws.Range("N92").Value = -19329.25
$ws.Range("H108").Value = 3017.3333
$ws.Range("J108").Value = 3017.3333
$ws.Range("L108").Value = 3017.3333
$ws.Range("N108").Value = -10697.3333
$ws.Range("H132").Value = 1197.2094
$ws.Range("I132").Value = 1053.0312
$ws.Range("K132").Value = 3159.0936
$ws.Range("M132").Value = -629.0935999999997
$ws.Range("H136").Value = 583.28
$ws.Range("I136").Value = 364.9
$ws.Range("J136").Value = 1456.8
$ws.Range("K136").Value = 1094.7
$ws.Range("L136").Value = 4370.4
$ws.Range("M136").Value = 1455.3
$ws.Range("N136").Value = -9470.4
